$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new expense entry in row 8 ("Resistors, rotary encoder, buttons" / Amazon / 49.35,
# dated 2016-03-30). The sheet has no data rows between row 7 and the TOTAL row at 21, so
# this just fills a previously-empty row rather than inserting/shifting anything.

# Column A (date) - copy the date format from the row above, then set the value.
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)
$ws.Cells.Item(8, 1).Value = Get-Date -Year 2016 -Month 3 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Column B (description)
$ws.Cells.Item(8, 2).Value = "Resistors, rotary encoder, buttons"

# Column C (vendor)
$ws.Cells.Item(8, 3).Value = "Amazon"

# Column D (price) - copy the currency format from the row above, then set the value.
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122)
$ws.Cells.Item(8, 4).Value = 49.35

$excel.CutCopyMode = $false

# The D21 TOTAL cell's SUM(D2:D20) formula already covers row 8, so it recalculates
# automatically and needs no explicit update.

# Move selection to B6, matching the authored change.
$ws.Range("B6").Select()
